# "so much joy today" -- refresh the RAND()-driven probability table on
# Sheet1 and tidy up the view the way the author left it.
#
# Column J (and therefore the dependent column K, including the running
# sums in K27/K40) is driven by volatile `RAND()*10` formulas, so simply
# forcing a recalculation gives each of those cells a brand new value --
# exactly what the commit's cell-value churn represents.  ironcalc's COM
# bridge already recalcs automatically once this script returns (it
# mirrors `Calculation = xlAutomatic`), but we also ask for it explicitly
# so the intent is obvious and the refresh happens even if that default
# ever changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force every volatile formula (RAND() in J2:J26, and the K-column/ totals
# that chain off of it) to recompute with fresh random draws.
$excel.CalculateFull()

# The author widened column J a bit (it holds the RAND()-based duration
# values) -- autofit it to its new best-fit width.
$ws.Columns.Item(10).EntireColumn.AutoFit()

# Scroll position reset to the top of the sheet and the selection moved
# to K2 (the first "weight" result cell) before the file was saved.
[void]$ws.Range("K2").Select()
